$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / title text updates (new report date 02-02-2025) ---
$ws.Range("A1").Value  = "Mangrove Communication   02.02.2025"
$ws.Range("A10").Value = "DAILY STOCK                         (02/02/2025) "

# --- Top sales table (rows 3-6) ---
$ws.Range("C3").Value = 28335
$ws.Range("D3").Value = 2070

$ws.Range("C4").Value = 29598
$ws.Range("D4").Value = 6963

$ws.Range("C5").Value = 30479
$ws.Range("D5").Value = 3011

$ws.Range("C6").Value = 28572
$ws.Range("D6").Value = 9411
$ws.Range("F6").Value = 10

# --- Stock table ---
# Row 14 - I top up
$ws.Range("C14").Value = 297170
$ws.Range("D14").Value = 116984
$ws.Range("E14").Value = 67974

# Row 20 - 19 tk mb
$ws.Range("C20").Value = 610
$ws.Range("D20").Value = 1140
$ws.Range("E20").Value = 2500

# Row 21 - 29 tk data (now cleared)
$ws.Range("C21").Value = ""
$ws.Range("D21").Value = ""

# Row 22 - 19 tk voice
$ws.Range("C22").Value = 450

# Row 24 - 50 tk
$ws.Range("C24").Value = 12

# Row 26 - Rbsp sim
$ws.Range("C26").Value = 28
$ws.Range("D26").Value = 10

# Row 27 - EV Swap Sim
$ws.Range("C27").Value = 83
$ws.Range("D27").Value = ""

# --- View state: keep selection in sync with the source workbook ---
$ws.Range("J21").Select()
